$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "67.098.30"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "3.128.18"
$ws.Range("E3").Value = "  +2.43%  "
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue $ws.Range("D5") "578.50"
$ws.Range("E5").Value = "  -0.04%  "
Set-TextValue $ws.Range("D6") "174.24"
$ws.Range("E6").Value = "  +3.75%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.123.57"
$ws.Range("E8").Value = "  +2.47%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E10").Value = "  -3.41%  "
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("E13").Value = "  +0.02%  "
Set-TextValue $ws.Range("D14") "37.27"
$ws.Range("E14").Value = "  +1.23%  "
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").Value = "3.647.09"
$ws.Range("E16").Value = "  +2.64%  "
$ws.Range("D17").Value = "67.119.10"
$ws.Range("E17").Value = "  +1.09%  "
Set-TextValue $ws.Range("D18") "7.14"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("D19").Value = "3.130.30"
$ws.Range("E19").Value = "  +2.73%  "
Set-TextValue $ws.Range("D20") "16.15"
$ws.Range("E20").Value = "  -2.90%  "
Set-TextValue $ws.Range("D21") "484.29"
$ws.Range("E21").Value = "  +3.82%  "
$ws.Range("E22").Value = "  +0.35%  "
Set-TextValue $ws.Range("D23") "7.71"
$ws.Range("E23").Value = "  +3.59%  "
$ws.Range("E26").Value = "  +1.44%  "
Set-TextValue $ws.Range("D27") "10.06"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  +0.00%  "
Set-TextValue $ws.Range("D29") "8.00"
$ws.Range("E29").Value = "  -2.37%  "
$ws.Range("E30").Value = "  -1.50%  "
$ws.Range("E31").Value = "  +0.58%  "
Set-TextValue $ws.Range("D32") "28.79"
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("D33").Value = "0.0₃0991"
$ws.Range("E33").Value = "  -1.98%  "
$ws.Range("E34").Value = "  -2.30%  "
$ws.Range("E35").Value = "  +0.12%  "
Set-TextValue $ws.Range("D36") "5.89"
$ws.Range("E36").Value = "  +0.21%  "
Set-TextValue $ws.Range("D37") "0.983"
$ws.Range("E37").Value = "  -1.30%  "
Set-TextValue $ws.Range("D38") "47.66"
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("E39").Value = "  +2.61%  "
Set-TextValue $ws.Range("D40") "50.12"
$ws.Range("E40").Value = "  +0.78%  "
Set-TextValue $ws.Range("D41") "0.312"
$ws.Range("E41").Value = "  -2.87%  "
Set-TextValue $ws.Range("D42") "0.123"
$ws.Range("E42").Value = "  +1.54%  "
$ws.Range("E43").Value = "  -0.40%  "
Set-TextValue $ws.Range("D44") "2.68"
$ws.Range("E44").Value = "  -5.61%  "
$ws.Range("D45").Value = "2.848.08"
$ws.Range("E45").Value = "  +3.91%  "
$ws.Range("E46").Value = "  -0.88%  "
Set-TextValue $ws.Range("D47") "382.59"
$ws.Range("E47").Value = "  -0.06%  "
Set-TextValue $ws.Range("D48") "136.05"
$ws.Range("E48").Value = "  +1.50%  "
$ws.Range("E49").Value = "  +0.01%  "
Set-TextValue $ws.Range("D50") "24.90"
$ws.Range("E50").Value = "  +1.39%  "
$ws.Range("E51").Value = "  -0.66%  "

# Rows 24/25 swap: Litecoin <-> InternetComputer(DFINITY)
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D24") "13.34"
$ws.Range("E24").Value = "  +3.60%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D25") "83.99"
$ws.Range("E25").Value = "  +0.97%  "
